$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table with updated
# attendance counts (column F) for rows 3-5.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1438
    $ws.Range("F4").Value = 95
    $ws.Range("F5").Value = 74
}
